# Rewrite the Bescheid body: subject changes from a Denkmalschutz repair order
# to a Gaststaettengesetz permit decision (new addressee, Tenor, Begruendung,
# Rechtsbehelfsbelehrung and signature block).
$d = $word.ActiveDocument

# Each element below becomes one `w:t` run, separated in the document by a
# blank line (two manual line breaks, PowerShell `` `v `` = Word's vertical-tab
# line-break char). Building the text this way (vs. Find&Replace) keeps the
# straight quotes around "Walfisch" from being auto-corrected into smart quotes.
$blocks = @(
    "Einleitung:",
    "Sehr geehrter Herr Graeter,",
    "Sie haben einen Antrag auf Erteilung einer Erlaubnis für eine Musikkneipe in den Räumen der ehemaligen Gaststätte `"Walfisch`" in Kehl gestellt. Nach Prüfung des Sachverhalts und der rechtlichen Voraussetzungen wird Ihnen hiermit die Erlaubnis unter bestimmten Bedingungen erteilt.",
    "Tenor:",
    "Die Erlaubnis für den Betrieb einer Musikkneipe in den Räumen der ehemaligen Gaststätte `"Walfisch`" in Kehl wird unter der Bedingung erteilt, dass die Herrentoilette entsprechend umgebaut wird. Bei Lärmbelästigungen in der Umgebung müssen Maßnahmen ergriffen werden.",
    "Begründung:",
    "Die Erlaubniserteilung erfolgt aufgrund der rechtlichen Voraussetzungen des § 12 GastG. Die Erlaubnispflicht ergibt sich aus § 2 GastG. Die Zuverlässigkeit des Antragstellers ist gemäß § 4 GastG geprüft worden und die Sachkunde gemäß § 5 GastG.",
    "Die Stadt Kehl hat gemäß § 12 GastG ein Ermessen, das nach § 40 LVwVfG ausgeübt wird. Die Erlaubniserteilung ist verhältnismäßig und ermessensgerecht, wenn die materiellen Voraussetzungen erfüllt sind und die öffentlichen Interessen nicht entgegenstehen.",
    "Es könnte eine Unmöglichkeit vorliegen, wenn die Herrentoilette nicht den Anforderungen entspricht. In diesem Fall könnte die Erlaubnis nur unter der Bedingung erteilt werden, dass die Toilette entsprechend umgebaut wird.",
    "Die Erlaubnis ist bestimmt genug formuliert, um den Anforderungen des § 37 LVwVfG zu genügen.",
    "Rechtsbehelfsbelehrung:",
    "Gegen diesen Bescheid kann innerhalb eines Monats nach Bekanntgabe Widerspruch eingelegt werden (§ 70 VwGO). Der Widerspruch ist schriftlich oder zur Niederschrift beim Amt für öffentliche Ordnung der Stadt Kehl, Hauptstraße 1, 77694 Kehl, einzulegen.",
    "Unterschrift mit Grußformel:",
    "Mit freundlichen Grüßen",
    "[Name und Funktion des Unterzeichners]"
)

$break = "`v`v"
$d.Content.Text = [string]::Join($break, $blocks)
